$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 85; this pushes the previous rows
# 85..124 down to 86..125 (the tail row 125 duplicates what used to be row 124).
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new record.
$ws.Cells.Item(85, 1).Value = 10
$ws.Cells.Item(85, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(85, 3).Value = "La Araucanía"
$ws.Cells.Item(85, 4).Value = 44813
$ws.Cells.Item(85, 5).Value = 9
$ws.Cells.Item(85, 6).Value = 100114002
$ws.Cells.Item(85, 7).Value = "Camote"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 30
$ws.Cells.Item(85, 11).Value = 20000
$ws.Cells.Item(85, 12).Value = 20000
$ws.Cells.Item(85, 13).Value = 20000
$ws.Cells.Item(85, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(85, 15).Value = "Perú"
$ws.Cells.Item(85, 16).Value = 1000
$ws.Cells.Item(85, 17).Value = 20
$ws.Cells.Item(85, 18).Value = "Hortaliza"
